$d = $word.ActiveDocument

# The resume's headings (name, "Berufserfahrung" section title, and each
# job/degree title) need to be made bold. Their paragraph-level style was
# already bold, but the run itself carried an explicit "not bold" override
# (<w:b w:val="0"/>) that suppressed it. Flip each of those runs to bold.
#
# "Senior Animation Designer" occurs twice: once as the subtitle under the
# name (paragraph 2, left untouched) and once as a job-title heading further
# down (paragraph 12, which must become bold). Match by paragraph index so
# only the intended heading runs are affected.
$targetParagraphIndexes = @(1, 5, 6, 12, 18, 24)

foreach ($i in $targetParagraphIndexes) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.Bold = 1
}
